$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("voc_site_classifications")

# Insert a new column before the existing "Classification" column (D),
# shifting it (and its data) one column to the right (D -> E).
$ws.Range("D1").EntireColumn.Insert()

# Header for the newly inserted column.
$ws.Range("D1").Value = "Classification Id"

# Assign a "CLn" identifier to each distinct Classification value (now in
# column E), numbered in the order each distinct value first appears.
$map = @{}
$nextNum = 1
for ($r = 2; $r -le 127; $r++) {
    $val = $ws.Cells.Item($r, 5).Value()
    if (-not $map.ContainsKey($val)) {
        $map[$val] = "CL" + $nextNum
        $nextNum = $nextNum + 1
    }
    $ws.Cells.Item($r, 4).Value = $map[$val]
}

# Keep the AutoFilter range in sync with the newly added column.
$ws.Range("A1:E1").AutoFilter()
